# SR [2022-09-02]: Some bugs fixed --> update Specification in existing files, if smth is changed in mngm xls
#
# Changes applied to "Planned Objects" sheet:
#  - J6: add a clarification comment
#  - Row 11 ("Some New Table" project row):
#      A11: "Some New Table" -> "Custom Internet Sales"
#      H11 (08_Status): "To Start" -> "Removed"
#      J11 (10_Dev Comment): clear the old "Waiting for specification" note
#      K11 (11_PM Comment): add removal note
#  - Move the active selection to J7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planned Objects")

$ws.Range("A11").Value = "Custom Internet Sales"
$ws.Range("H11").Value = "Removed"
$ws.Range("J11").ClearContents()
$ws.Range("K11").Value = "Removed from project on [2022-08-25] according to FUp letter from <CustomerRepresentative>"

$ws.Range("J6").Value = "Some clarification is required!"

$ws.Range("J7").Select()
